$d = $word.ActiveDocument

$pairs = @(
    @("56×32=1792", "20×46=920"),
    @("16×98=1568", "73×79=5767"),
    @("65×22=1430", "86×46=3956"),
    @("31×79=2449", "38×53=2014"),
    @("85×97=8245", "28×83=2324"),
    @("55×47=2585", "25×54=1350"),
    @("84×43=3612", "14×73=1022"),
    @("45×11=495",  "89×38=3382"),
    @("88×34=2992", "16×21=336"),
    @("30×37=1110", "44×49=2156"),
    @("37×68=2516", "63×76=4788"),
    @("23×17=391",  "33×23=759"),
    @("15×50=750",  "63×45=2835"),
    @("53×83=4399", "70×97=6790"),
    @("37×81=2997", "92×84=7728"),
    @("49×54=2646", "85×75=6375"),
    @("28×37=1036", "61×41=2501"),
    @("91×52=4732", "63×55=3465"),
    @("35×52=1820", "70×51=3570"),
    @("96×26=2496", "71×73=5183"),
    @("58×57=3306", "22×24=528"),
    @("88×55=4840", "33×87=2871"),
    @("39×87=3393", "34×86=2924"),
    @("20×32=640",  "30×98=2940"),
    @("66×21=1386", "13×64=832")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
